$wb = $excel.ActiveWorkbook

# Add "Dashboard" sheet right after the existing "Sheet"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$dashboard = $wb.Worksheets.Add($null, $lastSheet)
$dashboard.Name = "Dashboard"

# Add "Links" sheet right after "Dashboard"
$links = $wb.Worksheets.Add($null, $dashboard)
$links.Name = "Links"

# Populate the Dashboard sheet header
$dashboard.Range("A1").Value = "TotalNumber"

# Populate the Links sheet headers
$links.Range("A1").Value = "Link"
$links.Range("B1").Value = "DateAdded"
$links.Range("C1").Value = "Processed"
$links.Range("D1").Value = "DateProcessed"
